$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 174
$ws.Range("I5").Value = 58.75
$ws.Range("K5").Value = 58.75
$ws.Range("M5").Value = 56.25
$ws.Range("H127").Value = 790.7857
$ws.Range("I127").Value = 511.83334
$ws.Range("J127").Value = 1000
$ws.Range("K127").Value = 1535.50002
$ws.Range("L127").Value = 3000
$ws.Range("M127").Value = 3424.49998
$ws.Range("N127").Value = -12920
$ws.Range("H137").Value = 33334496
$ws.Range("I137").Value = 43479280
$ws.Range("J137").Value = 1629.4286
$ws.Range("K137").Value = 130437840
$ws.Range("L137").Value = 4888.2858
$ws.Range("M137").Value = -130435290
$ws.Range("N137").Value = -9988.2858
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 51419.15
$ws.Range("I2").Value = 68028
$ws.Range("J2").Value = 1592.6
$ws.Range("K2").Value = 68028
$ws.Range("L2").Value = 1592.6
$ws.Range("M2").Value = -67915
$ws.Range("N2").Value = -1818.6
$ws.Range("H23").Value = 17333.334
$ws.Range("J23").Value = 17333.334
$ws.Range("L23").Value = 17333.334
$ws.Range("N23").Value = -17851.334
$ws.Range("H63").Value = 9523.412
$ws.Range("I63").Value = 13126.546
$ws.Range("J63").Value = 2917.6667
$ws.Range("K63").Value = 13126.546
$ws.Range("L63").Value = 2917.6667
$ws.Range("M63").Value = -12440.546
$ws.Range("N63").Value = -4289.6667
$ws.Range("H66").Value = 9523.412
$ws.Range("I66").Value = 13126.546
$ws.Range("J66").Value = 2917.6667
$ws.Range("K66").Value = 65632.73
$ws.Range("L66").Value = 14588.3335
$ws.Range("M66").Value = -62200.73
$ws.Range("N66").Value = -21452.3335
$ws.Range("H74").Value = 4406.878
$ws.Range("I74").Value = 1368
$ws.Range("J74").Value = 11750.833
$ws.Range("K74").Value = 1368
$ws.Range("L74").Value = 11750.833
$ws.Range("M74").Value = -494
$ws.Range("N74").Value = -13498.833
$ws.Range("H77").Value = 4406.878
$ws.Range("I77").Value = 1368
$ws.Range("J77").Value = 11750.833
$ws.Range("K77").Value = 6840
$ws.Range("L77").Value = 58754.165
$ws.Range("M77").Value = -2472
$ws.Range("N77").Value = -67490.16500000001
$ws.Range("H116").Value = 51419.15
$ws.Range("I116").Value = 68028
$ws.Range("J116").Value = 1592.6
$ws.Range("K116").Value = 68028
$ws.Range("L116").Value = 1592.6
$ws.Range("M116").Value = -65734
$ws.Range("N116").Value = -6180.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 51419.15
$ws.Range("I3").Value = 68028
$ws.Range("J3").Value = 1592.6
$ws.Range("K3").Value = 68028
$ws.Range("L3").Value = 1592.6
$ws.Range("M3").Value = -67914
$ws.Range("N3").Value = -1820.6
$ws.Range("H94").Value = 1154.4
$ws.Range("I94").Value = 1104.5385
$ws.Range("J94").Value = 1247
$ws.Range("K94").Value = 1104.5385
$ws.Range("L94").Value = 1247
$ws.Range("M94").Value = -653.5385000000001
$ws.Range("N94").Value = -2149
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1570.35
$ws.Range("I31").Value = 994.6429000000001
$ws.Range("J31").Value = 2913.6667
$ws.Range("K31").Value = 994.6429000000001
$ws.Range("L31").Value = 2913.6667
$ws.Range("M31").Value = -699.6429000000001
$ws.Range("N31").Value = -3503.6667
$ws.Range("H34").Value = 1570.35
$ws.Range("I34").Value = 994.6429000000001
$ws.Range("J34").Value = 2913.6667
$ws.Range("K34").Value = 994.6429000000001
$ws.Range("L34").Value = 2913.6667
$ws.Range("M34").Value = -792.6429000000001
$ws.Range("N34").Value = -3317.6667
$ws.Range("H86").Value = 22729104
$ws.Range("I86").Value = 35716130
$ws.Range("J86").Value = 1812.25
$ws.Range("K86").Value = 35716130
$ws.Range("L86").Value = 1812.25
$ws.Range("M86").Value = -35715007
$ws.Range("N86").Value = -4058.25
$ws.Range("H89").Value = 22729104
$ws.Range("I89").Value = 35716130
$ws.Range("J89").Value = 1812.25
$ws.Range("K89").Value = 178580650
$ws.Range("L89").Value = 9061.25
$ws.Range("M89").Value = -178575034
$ws.Range("N89").Value = -20293.25
$ws.Range("H132").Value = 2111.3
$ws.Range("I132").Value = 1750.4375
$ws.Range("J132").Value = 3554.75
$ws.Range("K132").Value = 5251.3125
$ws.Range("L132").Value = 10664.25
$ws.Range("M132").Value = -2721.3125
$ws.Range("N132").Value = -15724.25
$ws.Range("H134").Value = 1936.2963
$ws.Range("I134").Value = 1229.475
$ws.Range("J134").Value = 3955.7856
$ws.Range("K134").Value = 3688.425
$ws.Range("L134").Value = 11867.3568
$ws.Range("M134").Value = -1153.425
$ws.Range("N134").Value = -16937.3568
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1180.4375
$ws.Range("I5").Value = 513.2222
$ws.Range("K5").Value = 1539.6666
$ws.Range("M5").Value = -1427.6666
$ws.Range("H135").Value = 1180.4375
$ws.Range("I135").Value = 513.2222
$ws.Range("K135").Value = 4618.999800000001
$ws.Range("M135").Value = -2083.999800000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1390588.9
$ws.Range("J122").Value = 1920
$ws.Range("L122").Value = 5760
$ws.Range("N122").Value = -10660
$ws.Range("H126").Value = 3255.55
$ws.Range("I126").Value = 3800
$ws.Range("J126").Value = 3074.0667
$ws.Range("K126").Value = 11400
$ws.Range("L126").Value = 9222.2001
$ws.Range("M126").Value = -8930
$ws.Range("N126").Value = -14162.2001
$ws.Range("H132").Value = 4773.579
$ws.Range("I132").Value = 4860.4287
$ws.Range("J132").Value = 4530.4
$ws.Range("K132").Value = 14581.2861
$ws.Range("L132").Value = 13591.2
$ws.Range("M132").Value = -12051.2861
$ws.Range("N132").Value = -18651.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4558.8237
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 4558.8237
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 4558.8237
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -4830.8237
$ws.Range("H136").Value = 4711.1055
$ws.Range("I136").Value = 2789.3794
$ws.Range("J136").Value = 10903.333
$ws.Range("K136").Value = 8368.138199999999
$ws.Range("L136").Value = 32709.999
$ws.Range("M136").Value = -5818.138199999999
$ws.Range("N136").Value = -37809.999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 802377.5
$ws.Range("I81").Value = 1178114
$ws.Range("J81").Value = 3937.5
$ws.Range("K81").Value = 2356228
$ws.Range("L81").Value = 7875
$ws.Range("M81").Value = -2355167
$ws.Range("N81").Value = -9997
$ws.Range("H84").Value = 802377.5
$ws.Range("I84").Value = 1178114
$ws.Range("J84").Value = 3937.5
$ws.Range("K84").Value = 11781140
$ws.Range("L84").Value = 39375
$ws.Range("M84").Value = -11775836
$ws.Range("N84").Value = -49983
$ws.Range("H122").Value = 57222.277
$ws.Range("I122").Value = 64225.062
$ws.Range("J122").Value = 1200
$ws.Range("K122").Value = 192675.186
$ws.Range("L122").Value = 3600
$ws.Range("M122").Value = -190225.186
$ws.Range("N122").Value = -8500
$ws.Range("H126").Value = 84719.5
$ws.Range("I126").Value = 200750.8
$ws.Range("J126").Value = 1840
$ws.Range("K126").Value = 602252.3999999999
$ws.Range("L126").Value = 5520
$ws.Range("M126").Value = -599782.3999999999
$ws.Range("N126").Value = -10460
$ws.Range("H132").Value = 17243846
$ws.Range("I132").Value = 20835666
$ws.Range("J132").Value = 3100.8
$ws.Range("K132").Value = 62506998
$ws.Range("L132").Value = 9302.400000000001
$ws.Range("M132").Value = -62504468
$ws.Range("N132").Value = -14362.4
$ws.Range("H136").Value = 9553971
$ws.Range("I136").Value = 12860083
$ws.Range("J136").Value = 2982.4443
$ws.Range("K136").Value = 38580249
$ws.Range("L136").Value = 8947.332900000001
$ws.Range("M136").Value = -38577699
$ws.Range("N136").Value = -14047.3329
